# Report template: clear the placeholder day-of-week numbers (1-7) that
# were left in the "Llamados" sub-header row, keeping their formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5:M5").ClearContents() | Out-Null

# Leave the cursor where the user last clicked before saving.
$ws.Range("M5").Select() | Out-Null
